# Rotate the data of rows 2, 3 and 4 (row 1 is the header row, untouched):
#   new row2 <- old row4
#   new row3 <- old row2
#   new row4 <- old row3
# Only the cells whose value actually differs between the old/new rows are
# written below (columns that are identical across all three source rows -
# C, D, N, P, S, T, U, V, W, Z, AB, AD, AE, AF, AG, AT, AW, AX, AY - are left
# untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 becomes old row 4 ----
$ws.Range("A2").Value = 104160759
$ws.Range("B2").Value = 90696
$ws.Range("E2").Value = 5448
$ws.Range("F2").Value = "Svartvit taggsvamp"
$ws.Range("G2").Value = "Phellodon connatus"
$ws.Range("H2").Value = "(Schultz) nom.prov"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "mycel"
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("Q2").Value = 549365.1106590288
$ws.Range("R2").Value = 6376508.764988168
$ws.Range("Y2").Value = "2022-10-16"
$ws.Range("AA2").Value = "2022-10-16"
$ws.Range("AC2").Value = ""

# ---- Row 3 becomes old row 2 ----
$ws.Range("A3").Value = 103437049
$ws.Range("B3").Value = 99566
$ws.Range("E3").Value = 221317
$ws.Range("F3").Value = "Gullklöver"
$ws.Range("G3").Value = "Trifolium aureum"
$ws.Range("H3").Value = "Pollich"
$ws.Range("I3").Value = ""
$ws.Range("K3").Value = ""
$ws.Range("M3").Value = ""
$ws.Range("Q3").Value = 549570.8693672653
$ws.Range("R3").Value = 6376356.619731385
$ws.Range("AC3").Value = "vägkant"

# ---- Row 4 becomes old row 3 ----
$ws.Range("A4").Value = 103437321
$ws.Range("B4").Value = 44332
$ws.Range("E4").Value = 102020
$ws.Range("F4").Value = "Smalsprötad bastardsvärmare"
$ws.Range("G4").Value = "Zygaena osterodensis"
$ws.Range("H4").Value = "Reiss, 1921"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "imago/adult"
$ws.Range("L4").Value = ""
$ws.Range("M4").Value = ""
$ws.Range("Q4").Value = 549410.0686067103
$ws.Range("R4").Value = 6376439.505282871
$ws.Range("Y4").Value = "2022-06-30"
$ws.Range("AA4").Value = "2022-06-30"
